$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 151.7260716666667
$ws.Range("H2").Value = 455.178215
$ws.Range("I2").Value = 0.2700739458961593
$ws.Range("J2").Value = 0.2783366498663096
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.670718
$ws.Range("N2").Value = 2.012154
$ws.Range("O2").Value = 0.05714337887504232
$ws.Range("P2").Value = 0.06502159750685171
$ws.Range("Q2").Value = 101.7654073361233
$ws.Range("R2").Value = 915.8886660251102
$ws.Range("S2").Value = 0.01543293781462191
$ws.Range("T2").Value = 0.01809789361901269

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 151.7260716666667
$ws.Range("H3").Value = 455.178215
$ws.Range("I3").Value = 0.2700739458961593
$ws.Range("J3").Value = 0.2783366498663096
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.04420533333333334
$ws.Range("N3").Value = 0.132616
$ws.Range("O3").Value = 0.003766176114200311
$ws.Range("P3").Value = 0.004285409653022904
$ws.Range("Q3").Value = 6.707101573382224
$ws.Range("R3").Value = 60.36391416044001
$ws.Range("S3").Value = 0.001017146044101942
$ws.Range("T3").Value = 0.001192786566127139

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 151.7260716666667
$ws.Range("H4").Value = 455.178215
$ws.Range("I4").Value = 0.2700739458961593
$ws.Range("J4").Value = 0.2783366498663096
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 3.258457666666667
$ws.Range("N4").Value = 9.775373
$ws.Range("O4").Value = 0.2776118741328243
$ws.Range("P4").Value = 0.3158855478682772
$ws.Range("Q4").Value = 494.3929814554662
$ws.Range("R4").Value = 4449.536833099195
$ws.Range("S4").Value = 0.07497573427467978
$ws.Range("T4").Value = 0.08792252513484004

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 151.7260716666667
$ws.Range("H5").Value = 455.178215
$ws.Range("I5").Value = 0.2700739458961593
$ws.Range("J5").Value = 0.2783366498663096
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 3.497636333333334
$ws.Range("N5").Value = 10.492909
$ws.Range("O5").Value = 0.2979892565322244
$ws.Range("P5").Value = 0.3390723104066696
$ws.Range("Q5").Value = 530.6826209752707
$ws.Range("R5").Value = 4776.143588777436
$ws.Range("S5").Value = 0.0804791343463207
$ws.Range("T5").Value = 0.09437625094102184

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 151.7260716666667
$ws.Range("H6").Value = 455.178215
$ws.Range("I6").Value = 0.2700739458961593
$ws.Range("J6").Value = 0.2783366498663096
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 4.2664405
$ws.Range("N6").Value = 8.532881
$ws.Range("O6").Value = 0.3634893143457086
$ws.Range("P6").Value = 0.2757351345651786
$ws.Range("Q6").Value = 647.3302570645692
$ws.Range("R6").Value = 3883.981542387415
$ws.Range("S6").Value = 0.09816899341643494
$ws.Range("T6").Value = 0.07674719360530786

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 82.24887099999999
$ws.Range("H7").Value = 246.746613
$ws.Range("I7").Value = 0.146403824289839
$ws.Range("J7").Value = 0.150882936320401
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.670718
$ws.Range("N7").Value = 2.012154
$ws.Range("O7").Value = 0.05714337887504232
$ws.Range("P7").Value = 0.06502159750685171
$ws.Range("Q7").Value = 55.165798259378
$ws.Range("R7").Value = 496.4921843344021
$ws.Range("S7").Value = 0.008366009200149392
$ws.Range("T7").Value = 0.00981064955607705

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 82.24887099999999
$ws.Range("H8").Value = 246.746613
$ws.Range("I8").Value = 0.146403824289839
$ws.Range("J8").Value = 0.150882936320401
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.04420533333333334
$ws.Range("N8").Value = 0.132616
$ws.Range("O8").Value = 0.003766176114200311
$ws.Range("P8").Value = 0.004285409653022904
$ws.Range("Q8").Value = 3.635838758845333
$ws.Range("R8").Value = 32.722548829608
$ws.Range("S8").Value = 0.0005513825860679708
$ws.Range("T8").Value = 0.0006465951917838864

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 82.24887099999999
$ws.Range("H9").Value = 246.746613
$ws.Range("I9").Value = 0.146403824289839
$ws.Range("J9").Value = 0.150882936320401
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 3.258457666666667
$ws.Range("N9").Value = 9.775373
$ws.Range("O9").Value = 0.2776118741328243
$ws.Range("P9").Value = 0.3158855478682772
$ws.Range("Q9").Value = 268.0044642846277
$ws.Range("R9").Value = 2412.040178561649
$ws.Range("S9").Value = 0.04064344004131491
$ws.Range("T9").Value = 0.04766173900354425

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 82.24887099999999
$ws.Range("H10").Value = 246.746613
$ws.Range("I10").Value = 0.146403824289839
$ws.Range("J10").Value = 0.150882936320401
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 3.497636333333334
$ws.Range("N10").Value = 10.492909
$ws.Range("O10").Value = 0.2979892565322244
$ws.Range("P10").Value = 0.3390723104066696
$ws.Range("Q10").Value = 287.6766395852463
$ws.Range("R10").Value = 2589.089756267217
$ws.Range("S10").Value = 0.04362676675360353
$ws.Range("T10").Value = 0.05116022581910077

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 82.24887099999999
$ws.Range("H11").Value = 246.746613
$ws.Range("I11").Value = 0.146403824289839
$ws.Range("J11").Value = 0.150882936320401
$ws.Range("K11").Value = 2
$ws.Range("M11").Value = 4.2664405
$ws.Range("N11").Value = 8.532881
$ws.Range("O11").Value = 0.3634893143457086
$ws.Range("P11").Value = 0.2757351345651786
$ws.Range("Q11").Value = 350.9099143136755
$ws.Range("R11").Value = 2105.459485882053
$ws.Range("S11").Value = 0.05321622570870316
$ws.Range("T11").Value = 0.04160372674989503

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 123.444321
$ws.Range("H12").Value = 370.332963
$ws.Range("I12").Value = 0.2197321429647646
$ws.Range("J12").Value = 0.2264546783208506
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.670718
$ws.Range("N12").Value = 2.012154
$ws.Range("O12").Value = 0.05714337887504232
$ws.Range("P12").Value = 0.06502159750685171
$ws.Range("Q12").Value = 82.79632809247801
$ws.Range("R12").Value = 745.1669528323021
$ws.Range("S12").Value = 0.01255623709646051
$ws.Range("T12").Value = 0.01472444494732193

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 123.444321
$ws.Range("H13").Value = 370.332963
$ws.Range("I13").Value = 0.2197321429647646
$ws.Range("J13").Value = 0.2264546783208506
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 0.04420533333333334
$ws.Range("N13").Value = 0.132616
$ws.Range("O13").Value = 0.003766176114200311
$ws.Range("P13").Value = 0.004285409653022904
$ws.Range("Q13").Value = 5.456897357912001
$ws.Range("R13").Value = 49.11207622120801
$ws.Range("S13").Value = 0.0008275499483559444
$ws.Range("T13").Value = 0.0009704510644483696

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 123.444321
$ws.Range("H14").Value = 370.332963
$ws.Range("I14").Value = 0.2197321429647646
$ws.Range("J14").Value = 0.2264546783208506
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 3.258457666666667
$ws.Range("N14").Value = 9.775373
$ws.Range("O14").Value = 0.2776118741328243
$ws.Range("P14").Value = 0.3158855478682772
$ws.Range("Q14").Value = 402.238094168911
$ws.Range("R14").Value = 3620.142847520199
$ws.Range("S14").Value = 0.06100025201566999
$ws.Range("T14").Value = 0.07153376012871637

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 123.444321
$ws.Range("H15").Value = 370.332963
$ws.Range("I15").Value = 0.2197321429647646
$ws.Range("J15").Value = 0.2264546783208506
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 3.497636333333334
$ws.Range("N15").Value = 10.492909
$ws.Range("O15").Value = 0.2979892565322244
$ws.Range("P15").Value = 0.3390723104066696
$ws.Range("Q15").Value = 431.763342273263
$ws.Range("R15").Value = 3885.870080459367
$ws.Range("S15").Value = 0.06547781791830265
$ws.Range("T15").Value = 0.07678451098064996

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 123.444321
$ws.Range("H16").Value = 370.332963
$ws.Range("I16").Value = 0.2197321429647646
$ws.Range("J16").Value = 0.2264546783208506
$ws.Range("K16").Value = 2
$ws.Range("M16").Value = 4.2664405
$ws.Range("N16").Value = 8.532881
$ws.Range("O16").Value = 0.3634893143457086
$ws.Range("P16").Value = 0.2757351345651786
$ws.Range("Q16").Value = 526.6678506094005
$ws.Range("R16").Value = 3160.007103656403
$ws.Range("S16").Value = 0.07987028598597549
$ws.Range("T16").Value = 0.06244151119971396

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 154.3429766666667
$ws.Range("H17").Value = 463.02893
$ws.Range("I17").Value = 0.2747320633285943
$ws.Range("J17").Value = 0.2831372788071194
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 0.670718
$ws.Range("N17").Value = 2.012154
$ws.Range("O17").Value = 0.05714337887504232
$ws.Range("P17").Value = 0.06502159750685171
$ws.Range("Q17").Value = 103.5206126239133
$ws.Range("R17").Value = 931.6855136152201
$ws.Range("S17").Value = 0.01569911838390798
$ws.Range("T17").Value = 0.01841003818178177

# Row 18
$ws.Range("E18").Value = 3
$ws.Range("G18").Value = 154.3429766666667
$ws.Range("H18").Value = 463.02893
$ws.Range("I18").Value = 0.2747320633285943
$ws.Range("J18").Value = 0.2831372788071194
$ws.Range("K18").Value = 3
$ws.Range("M18").Value = 0.04420533333333334
$ws.Range("N18").Value = 0.132616
$ws.Range("O18").Value = 0.003766176114200311
$ws.Range("P18").Value = 0.004285409653022904
$ws.Range("Q18").Value = 6.822782731208889
$ws.Range("R18").Value = 61.40504458088001
$ws.Range("S18").Value = 0.001034689334713119
$ws.Range("T18").Value = 0.001213359227730667

# Row 19
$ws.Range("E19").Value = 3
$ws.Range("G19").Value = 154.3429766666667
$ws.Range("H19").Value = 463.02893
$ws.Range("I19").Value = 0.2747320633285943
$ws.Range("J19").Value = 0.2831372788071194
$ws.Range("K19").Value = 3
$ws.Range("M19").Value = 3.258457666666667
$ws.Range("N19").Value = 9.775373
$ws.Range("O19").Value = 0.2776118741328243
$ws.Range("P19").Value = 0.3158855478682772
$ws.Range("Q19").Value = 502.9200556156544
$ws.Range("R19").Value = 4526.28050054089
$ws.Range("S19").Value = 0.07626888298502885
$ws.Range("T19").Value = 0.08943897443792008

# Row 20
$ws.Range("E20").Value = 3
$ws.Range("G20").Value = 154.3429766666667
$ws.Range("H20").Value = 463.02893
$ws.Range("I20").Value = 0.2747320633285943
$ws.Range("J20").Value = 0.2831372788071194
$ws.Range("K20").Value = 3
$ws.Range("M20").Value = 3.497636333333334
$ws.Range("N20").Value = 10.492909
$ws.Range("O20").Value = 0.2979892565322244
$ws.Range("P20").Value = 0.3390723104066696
$ws.Range("Q20").Value = 539.8356029841523
$ws.Range("R20").Value = 4858.520426857371
$ws.Range("S20").Value = 0.08186720329685181
$ws.Range("T20").Value = 0.09600401128738736

# Row 21
$ws.Range("E21").Value = 3
$ws.Range("G21").Value = 154.3429766666667
$ws.Range("H21").Value = 463.02893
$ws.Range("I21").Value = 0.2747320633285943
$ws.Range("J21").Value = 0.2831372788071194
$ws.Range("K21").Value = 2
$ws.Range("M21").Value = 4.2664405
$ws.Range("N21").Value = 8.532881
$ws.Range("O21").Value = 0.3634893143457086
$ws.Range("P21").Value = 0.2757351345651786
$ws.Range("Q21").Value = 658.4951265412216
$ws.Range("R21").Value = 3950.97075924733
$ws.Range("S21").Value = 0.09986216932809253
$ws.Range("T21").Value = 0.07807089567229955

# Row 22
$ws.Range("E22").Value = 2
$ws.Range("G22").Value = 50.0323125
$ws.Range("H22").Value = 100.064625
$ws.Range("I22").Value = 0.08905802352064279
$ws.Range("J22").Value = 0.06118845668531954
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 0.670718
$ws.Range("N22").Value = 2.012154
$ws.Range("O22").Value = 0.05714337887504232
$ws.Range("P22").Value = 0.06502159750685171
$ws.Range("Q22").Value = 33.557572575375
$ws.Range("R22").Value = 201.34543545225
$ws.Range("S22").Value = 0.005089076379902521
$ws.Range("T22").Value = 0.003978571202658277

# Row 23
$ws.Range("E23").Value = 2
$ws.Range("G23").Value = 50.0323125
$ws.Range("H23").Value = 100.064625
$ws.Range("I23").Value = 0.08905802352064279
$ws.Range("J23").Value = 0.06118845668531954
$ws.Range("K23").Value = 3
$ws.Range("M23").Value = 0.04420533333333334
$ws.Range("N23").Value = 0.132616
$ws.Range("O23").Value = 0.003766176114200311
$ws.Range("P23").Value = 0.004285409653022904
$ws.Range("Q23").Value = 2.2116950515
$ws.Range("R23").Value = 13.270170309
$ws.Range("S23").Value = 0.0003354082009613344
$ws.Range("T23").Value = 0.0002622176029328422

# Row 24
$ws.Range("E24").Value = 2
$ws.Range("G24").Value = 50.0323125
$ws.Range("H24").Value = 100.064625
$ws.Range("I24").Value = 0.08905802352064279
$ws.Range("J24").Value = 0.06118845668531954
$ws.Range("K24").Value = 3
$ws.Range("M24").Value = 3.258457666666667
$ws.Range("N24").Value = 9.775373
$ws.Range("O24").Value = 0.2776118741328243
$ws.Range("P24").Value = 0.3158855478682772
$ws.Range("Q24").Value = 163.0281722466875
$ws.Range("R24").Value = 978.1690334801251
$ws.Range("S24").Value = 0.0247235648161308
$ws.Range("T24").Value = 0.01932854916325651

# Row 25
$ws.Range("E25").Value = 2
$ws.Range("G25").Value = 50.0323125
$ws.Range("H25").Value = 100.064625
$ws.Range("I25").Value = 0.08905802352064279
$ws.Range("J25").Value = 0.06118845668531954
$ws.Range("K25").Value = 3
$ws.Range("M25").Value = 3.497636333333334
$ws.Range("N25").Value = 10.492909
$ws.Range("O25").Value = 0.2979892565322244
$ws.Range("P25").Value = 0.3390723104066696
$ws.Range("Q25").Value = 174.9948340406875
$ws.Range("R25").Value = 1049.969004244125
$ws.Range("S25").Value = 0.0265383342171457
$ws.Range("T25").Value = 0.02074731137850973

# Row 26
$ws.Range("E26").Value = 2
$ws.Range("G26").Value = 50.0323125
$ws.Range("H26").Value = 100.064625
$ws.Range("I26").Value = 0.08905802352064279
$ws.Range("J26").Value = 0.06118845668531954
$ws.Range("K26").Value = 2
$ws.Range("M26").Value = 4.2664405
$ws.Range("N26").Value = 8.532881
$ws.Range("O26").Value = 0.3634893143457086
$ws.Range("P26").Value = 0.2757351345651786
$ws.Range("Q26").Value = 213.4598843586563
$ws.Range("R26").Value = 853.8395374346251
$ws.Range("S26").Value = 0.03237163990650244
$ws.Range("T26").Value = 0.01687180733796218
